$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fromCSV")

# B80: "PsE6uD" -> "3Xx1rt"
$ws.Range("B80").Value = "3Xx1rt"

# O80: "106" -> "105"
$ws.Range("O80").Value = "105"
